# login page code refactoring
# Re-order the login-attempt rows (2-5), update the Pass/Fail outcomes,
# drop the hyperlink from two of the now-non-linked rows, and give row 2's
# "A" cell its own (smaller-font) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the row values (rows 2-5 of columns A:C) -------------------
# New layout:
#   Row2: admin@yourstore.com | adm   | Fail
#   Row3: adm@yourstore.com   | admin | Fail
#   Row4: adm@yourstore.com   | adm   | Fail
#   Row5: admin@yourstore.com | admin | Pass

$ws.Range("A2").Value = "admin@yourstore.com "
$ws.Range("B2").Value = "adm"
$ws.Range("C2").Value = "Fail"

$ws.Range("A3").Value = "adm@yourstore.com"
$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "Fail"

$ws.Range("A4").Value = "adm@yourstore.com"
$ws.Range("B4").Value = "adm"
$ws.Range("C4").Value = "Fail"

$ws.Range("A5").Value = "admin@yourstore.com "
$ws.Range("B5").Value = "admin"
$ws.Range("C5").Value = "Pass"

# --- 2. Give A2 its own style (smaller hyperlink-weight font) -------------
# It keeps the centered alignment + border that the whole column already
# uses; only the font shrinks from the big (size 14) hyperlink font to the
# regular (size 11) hyperlink font.
$ws.Range("A2").Font.Size = 11

# --- 3. Hyperlinks: A3 and A4 are no longer links; A2/A5 keep theirs ------
# (A2 already pointed at admin@yourstore.com and A5 already pointed at
# adm@yourstore.com, so only the now-orphaned links need to be removed.)
$toRemove = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 3 -or $h.Range.Row -eq 4) {
        $toRemove += $h
    }
}
foreach ($h in $toRemove) {
    $h.Delete()
}

# --- 4. Misc view state ----------------------------------------------------
$ws.Range("C3").Select()
